$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$text)
    # Force the cell to be stored as text even when the string looks numeric
    # (mirrors the original file, which keeps values like "780" as text).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Updated / new rows: row, A (Residue Combination), B (Count), C (Permeated Residues), D (Frames)
$data = @(
    @(2,  "98, 130, 748, 1073",        3, "130, 130, 130",   "5178, 5400, 5582"),
    @(3,  "98, 455, 780, 1105",        1, "780",              "5552"),
    @(4,  "130, 780, 780",             1, "130",              "3171"),
    @(5,  "98, 130, 748, 780",         1, "130",              "3631"),
    @(6,  "98, 780, 1105",             1, "780",              "3666"),
    @(7,  "130, 455, 748, 1073",       1, "130",              "5131"),
    @(8,  "98, 130, 455, 1073",        2, "130, 130",         "4416, 6489"),
    @(9,  "130, 423, 1073, SF, SF",    1, "130",              "4995"),
    @(10, "130, 455, 748, 780",        1, "130",              "5269"),
    @(11, "423, 748, 780, 1073",       1, "780",              "5677"),
    @(12, "423, 748, 1073, 1105",      2, "1105, 1105",       "5331, 5433"),
    @(13, "130, 423, 748, 1073",       1, "130",              "5886"),
    @(14, "98, 130, 748, 780, 1073",   1, "130",              "6017"),
    @(15, "130, 423, 748, 1073, SF",   1, "130",              "6202"),
    @(16, "98, 130, 423",              1, "130",              "6427"),
    @(17, "130, 455, 780, 1073",       1, "130",              "6562"),
    @(18, "98, 780, 780, 1105",        1, "780",              "6359"),
    @(19, "98, 130, 130, 455, 780",    1, "130",              "6727"),
    @(20, "130, 130, 423, 1073",       1, "130",              "6670"),
    @(21, "98, 98, 455, 455",          1, "455",              "6748")
)

foreach ($row in $data) {
    $r = $row[0]

    Set-TextCell $ws.Cells.Item($r, 1) $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    Set-TextCell $ws.Cells.Item($r, 3) $row[3]
    Set-TextCell $ws.Cells.Item($r, 4) $row[4]
}
